$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the now-empty cells (rows 2-6 in column C and E)
$ws.Range("C2:C6").ClearContents()
$ws.Range("E2:E6").ClearContents()

# Update the C column (y_0_forecast) values for rows 7-19
$ws.Range("C7").Value = 4.880442637054072
$ws.Range("C8").Value = 5.941867202078877
$ws.Range("C9").Value = 0.292749233164491
$ws.Range("C10").Value = 2.032207428223742
$ws.Range("C11").Value = 2.775332754349846
$ws.Range("C12").Value = 3.565025829754953
$ws.Range("C13").Value = 3.444206290325491
$ws.Range("C14").Value = 2.667234932970275
$ws.Range("C15").Value = -4.511102905979703
$ws.Range("C16").Value = 1.386772772629241
$ws.Range("C17").Value = -0.9537175292835154
$ws.Range("C18").Value = -3.303819519576723
$ws.Range("C19").Value = -1.704805397136089

# Update the E column (y_1_forecast) values for rows 7-19
$ws.Range("E7").Value = 1.641301872652501
$ws.Range("E8").Value = 2.672847571394987
$ws.Range("E9").Value = 2.181874035977249
$ws.Range("E10").Value = 2.13692496326825
$ws.Range("E11").Value = 2.200426660963761
$ws.Range("E12").Value = 2.446228176258058
$ws.Range("E13").Value = 2.496958452261078
$ws.Range("E14").Value = 2.567662999186382
$ws.Range("E15").Value = 0.9582724917052587
$ws.Range("E16").Value = 0.8813242377093244
$ws.Range("E17").Value = 0.9049225073274991
$ws.Range("E18").Value = -0.3655818470008065
$ws.Range("E19").Value = 0.1460701281005727
